# MITHEx_inputs.xlsx update
# - Remove excess print statements (no spreadsheet-visible effect)
# - Update pump and turbine equations (Compression Ratio parameter change)
# - Allow for use of Supercritical CO2 (rename "SCO2" option -> "CarbonDioxide",
#   update default Secondary Fluid selection, update Channel Diameter default)
# - Add pressure drop outputs (Secondary Pressure bound increased)

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet: Input options
# Rename the "SCO2" fluid option to "CarbonDioxide"
# ---------------------------------------------------------------------------
$wsOptions = $wb.Worksheets.Item("Input options")
$wsOptions.Range("B3").Value = "CarbonDioxide"

# ---------------------------------------------------------------------------
# Sheet: Plant Description
# Update the Secondary Fluid selection and Secondary Pressure bound
# ---------------------------------------------------------------------------
$wsPlant = $wb.Worksheets.Item("Plant Description")
$wsPlant.Range("B7").Value = "CarbonDioxide"
$wsPlant.Range("B11").Value = 25000

# ---------------------------------------------------------------------------
# Sheet: HX Parameters
# Update the Channel Diameter default value
# ---------------------------------------------------------------------------
$wsHX = $wb.Worksheets.Item("HX Parameters")
$wsHX.Range("B2").Value = 0.001

# ---------------------------------------------------------------------------
# Sheet: Cycle Parameters
# Update the Compression Ratio default value
# ---------------------------------------------------------------------------
$wsCycle = $wb.Worksheets.Item("Cycle Parameters")
$wsCycle.Range("B3").Value = 3

# ---------------------------------------------------------------------------
# Restore / update the view state (selected cell on each sheet, and which
# sheet tab is active) to match the saved workbook state.
# ---------------------------------------------------------------------------
$wsPlant.Activate()
$wsPlant.Range("B12").Select()

$wsHX.Activate()
$wsHX.Range("B2").Select()

$wsOptions.Activate()
$wsOptions.Range("B4").Select()

$wsCycle.Activate()
$wsCycle.Range("B4").Select()

$win = $excel.ActiveWindow
$win.Top = 500
$win.Height = 20680
